$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "68.642.68"
$ws.Range("E2").Value = "  +2.28%  "
Set-TextValue "D3" "2.527.74"
$ws.Range("E3").Value = "  +2.39%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  +1.99%  "
Set-TextValue "D6" "177.22"
$ws.Range("E6").Value = "  +1.44%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +1.66%  "
Set-TextValue "D9" "2.527.09"
$ws.Range("E9").Value = "  +2.39%  "
Set-TextValue "D10" "0.145"
$ws.Range("E10").Value = "  +5.48%  "
Set-TextValue "D12" "4.99"
$ws.Range("E12").Value = "  +0.94%  "
Set-TextValue "D13" "0.339"
$ws.Range("E13").Value = "  +1.63%  "
Set-TextValue "D14" "2.989.28"
$ws.Range("E14").Value = "  +2.48%  "
Set-TextValue "D15" "26.22"
$ws.Range("E15").Value = "  +3.22%  "
Set-TextValue "D16" "68.581.16"
$ws.Range("E16").Value = "  +2.47%  "
Set-TextValue "D17" "0.0000170"
$ws.Range("E17").Value = "  +0.66%  "
Set-TextValue "D18" "2.534.18"
$ws.Range("E18").Value = "  +1.31%  "
Set-TextValue "D19" "11.08"
$ws.Range("E19").Value = "  +1.41%  "
Set-TextValue "D20" "7.51"
$ws.Range("E20").Value = "  +0.68%  "
Set-TextValue "D21" "352.22"
$ws.Range("E21").Value = "  +1.09%  "
Set-TextValue "D22" "4.21"
$ws.Range("E22").Value = "  +5.72%  "
Set-TextValue "D23" "0.999"
$ws.Range("E23").Value = "  -0.11%  "
Set-TextValue "D24" "70.93"
$ws.Range("E24").Value = "  +2.32%  "
$ws.Range("E25").Value = "  +1.43%  "
Set-TextValue "D26" "1.70"
$ws.Range("E26").Value = "  -5.45%  "
Set-TextValue "D27" "9.01"
$ws.Range("E27").Value = "  -2.41%  "
Set-TextValue "D28" "2.690.82"
$ws.Range("E28").Value = "  +3.69%  "
Set-TextValue "D29" "0.990"
$ws.Range("E29").Value = "  -0.97%  "
Set-TextValue "D30" "510.35"
$ws.Range("E30").Value = "  +2.30%  "
Set-TextValue "D31" "0.0₃0893"
$ws.Range("E31").Value = "  -0.92%  "
Set-TextValue "D32" "7.80"
$ws.Range("E32").Value = "  +0.89%  "
Set-TextValue "D33" "1.25"
$ws.Range("E33").Value = "  +1.82%  "
$ws.Range("E34").Value = "  +1.19%  "
Set-TextValue "D35" "0.999"
$ws.Range("E35").Value = "  -0.08%  "
Set-TextValue "D36" "162.95"
$ws.Range("E36").Value = "  +0.95%  "
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("B38").Value = "WhiteBITCoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D38" "18.68"
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D39" "18.40"
$ws.Range("E39").Value = "  +1.34%  "
$ws.Range("E40").Value = "  +5.17%  "
$ws.Range("E41").Value = "  -0.47%  "
$ws.Range("E42").Value = "  +0.00%  "
Set-TextValue "D43" "4.84"
$ws.Range("E43").Value = "  +0.40%  "
Set-TextValue "D44" "0.326"
$ws.Range("E44").Value = "  -0.29%  "
Set-TextValue "D45" "2.41"
$ws.Range("E45").Value = "  +0.82%  "
Set-TextValue "D46" "153.48"
$ws.Range("E46").Value = "  +7.66%  "
$ws.Range("E47").Value = "  +2.65%  "
Set-TextValue "D48" "0.521"
$ws.Range("E48").Value = "  +2.27%  "
Set-TextValue "D49" "0.0₆0259"
$ws.Range("E49").Value = "  +1.22%  "
Set-TextValue "D50" "1.61"
$ws.Range("E50").Value = "  +2.59%  "
Set-TextValue "D51" "0.0740"
$ws.Range("E51").Value = "  -0.17%  "
